$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2, 3 and 5 of the sheet are cyclically rotated:
#   original row 5 data -> new row 2
#   original row 2 data -> new row 3
#   original row 3 data -> new row 5
# Row 4 is left completely untouched.
#
# Read every source value up-front (before any writes happen) so that the
# later in-place updates don't clobber data that still needs to be copied
# elsewhere.

$cols = @("A","B","D","E","F","G","H","I","Q","R","S","Z","AB")

$row2 = @{}
$row3 = @{}
$row5 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range("${col}2").Value()
    $row3[$col] = $ws.Range("${col}3").Value()
    $row5[$col] = $ws.Range("${col}5").Value()
}

# Only touch a cell when its value actually changes. This avoids rewriting
# cells that keep the same content across the rotation (e.g. D2 stays "NT",
# I2 stays blank, S5 stays 10), matching the source diff exactly.
foreach ($col in $cols) {
    if ($row5[$col] -ne $row2[$col]) {
        $ws.Range("${col}2").Value = $row5[$col]
    }
}
foreach ($col in $cols) {
    if ($row2[$col] -ne $row3[$col]) {
        $ws.Range("${col}3").Value = $row2[$col]
    }
}
foreach ($col in $cols) {
    if ($row3[$col] -ne $row5[$col]) {
        $ws.Range("${col}5").Value = $row3[$col]
    }
}
